# Add duplicate detection for contract note imports.
#
# A newly-imported contract note (CN#252611730667, dated one day after the
# existing entry) turned out to be a distinct trade rather than a duplicate
# of the existing row, so it is inserted above the existing "Trading
# History" entry that previously sat in row 5 (CN#252611665409). That
# existing entry moves down to row 6 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Move the existing row 5 entry (NSE / Buy / CN#252611665409, dated 46062)
# down to row 6 to make room for the newly detected entry above it.
$ws.Range("A6").Value2 = 46062
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat
$ws.Range("B6").Value2 = "NSE"
$ws.Range("C6").Value2 = "Buy"
$ws.Range("D6").Value2 = 5
$ws.Range("E6").Value2 = 891.95
$ws.Range("F6").Value2 = 4491.35
$ws.Range("G6").Value2 = "CN#252611665409"
$ws.Range("H6").Value2 = 4.46
$ws.Range("I6").Value2 = 27.14
$ws.Range("J6").Formula = "=Index!`$C`$2"

# Write the newly detected (non-duplicate) contract note entry into row 5.
$ws.Range("A5").Value2 = 46063
$ws.Range("B5").Value2 = "NSE"
$ws.Range("C5").Value2 = "Buy"
$ws.Range("D5").Value2 = 5
$ws.Range("E5").Value2 = 890.01
$ws.Range("F5").Value2 = 4481.6
$ws.Range("G5").Value2 = "CN#252611730667"
$ws.Range("H5").Value2 = 4.45
$ws.Range("I5").Value2 = 27.1
$ws.Range("J5").Formula = "=Index!`$C`$2"
